# "add next day to data" - a new day's row was appended upstream and the
# ticker workbook was refreshed. In the sheet itself this shows up as the
# color-filter on column B (which had been hiding every row except the
# "highlighted" tickers) being cleared, so all 122 previously-filtered rows
# (9-137) become visible again, and the summary formulas in rows 1-6 (which
# use SUBTOTAL over the visible rows) recompute against the full data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear every active AutoFilter criterion on the sheet (the color filter on
# column B in this case) and unhide the rows it was hiding. This mirrors a
# user clicking "Clear Filter" / "Select All" on the AutoFilter drop-down.
$ws.ShowAllData()

# Move the selection to where the author left off after reviewing the newly
# revealed data.
$ws.Range("CM50").Select()
